$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "row label" cell (A7, bold/centered/bordered)
# down onto the new row labels A8:A14 so no new cell styles get minted.
$ws.Range("A7").Copy()
$ws.Range("A8:A14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$data = @(
    @(8,  0.7489285731928166, 0.8279462115860953, 0.8055372423794177, 0.9906730972243635, 1.197432579484514,  0.7602408221395534),
    @(9,  1.10010523950836,   1.12279591879383,   1.537163856773588,  1.264524664494019,  1.1351442256939,    0.946760951598782),
    @(10, 0.8231969701403786, 1.164292629627047,  1.063463882421663,  1.101371799224184,  1.301639401483388,  1.06516361803197),
    @(11, 0.7308399557844771, 0.6789794128506424, 0.6116870057193515, 0.6490724813146383, 0.9481986454493596, 1.234403594205529),
    @(12, 0.9303135594904726, 1.174207595199213,  0.9822517680118442, 1.135300177435434,  1.281085910869446,  1.341625807083262),
    @(13, 0.898281320836375,  0.8523280122664625, 1.147700678224136,  1.050865029082537,  1.142395796240432,  1.074238764262095),
    @(14, 1.264939569612696,  1.11724155029655,   1.081276692631159,  1.443104718815984,  1.382728978810254,  1.360412685622154)
)

$rowIndex = 8
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($rowIndex, $col + 1).Value = $row[$col]
    }
    $rowIndex++
}
